$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 09:22"

# 2. Polonia (row 30): update active cases / recovered
$ws.Range("D30").Value = 618
$ws.Range("E30").Value = 6071

# 3. Moldavia (row 60): update active cases / recovered / deaths today / deaths
$ws.Range("D60").Value = 134
$ws.Range("E60").Value = 1542
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 36

# 4. Barein (row 66): update recovered / critical / deaths today / deaths
$ws.Range("E66").Value = 763
$ws.Range("F66").Value = 3
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 7

# 5. Rows 72-75 get re-sorted (by total cases, descending) and refreshed with
#    new daily figures for Bosnia y Herzegovina, Lituania, Armenia and Uzbekistan.
$ws.Range("A72").Value = "Bosnia y Herzegovina"
$ws.Range("B72").Value = 1077
$ws.Range("C72").Value = 40
$ws.Range("D72").Value = 218
$ws.Range("E72").Value = 820
$ws.Range("F72").Value = 4
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 39

$ws.Range("A73").Value = "Lituania"
$ws.Range("B73").Value = 1070
$ws.Range("C73").Value = 8
$ws.Range("D73").Value = 101
$ws.Range("E73").Value = 945
$ws.Range("F73").Value = 14
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 24

$ws.Range("A74").Value = "Armenia"
$ws.Range("B74").Value = 1067
$ws.Range("C74").Value = 28
$ws.Range("D74").Value = 265
$ws.Range("E74").Value = 786
$ws.Range("F74").Value = 30
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 16

$ws.Range("A75").Value = "Uzbekistan"
$ws.Range("B75").Value = 1054
$ws.Range("C75").Value = 56
$ws.Range("D75").Value = 89
$ws.Range("E75").Value = 961
$ws.Range("F75").Value = 8
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 4

# 6. Letonia (row 88): update total cases / new cases / recovered / critical
$ws.Range("B88").Value = 657
$ws.Range("C88").Value = 2
$ws.Range("E88").Value = 636
$ws.Range("F88").Value = 3
